$wb = $excel.ActiveWorkbook

# --- Sheets ---
$ws1 = $wb.Worksheets.Item("Paises")
$ws2 = $wb.Worksheets.Item("Clientes")
$ws3 = $wb.Worksheets.Item("ZeroCode")

# ======================================================================
# Sheet "Clientes": insert the new "Pagina web" row, reorder the
# "Observaciones"/"Fecha registro" rows, add "Hora registro" and
# "Ultima actualizacion" rows, and push "Recibir publicidad" to the end
# (clearing its "obligatorio" flag so it's blank by default).
# ======================================================================
$ws2.Rows("7:7").Insert()
$ws2.Rows("11:12").Insert()
$ws2.Range("A7:K13").ClearContents()

$ws2.Range("A7").Value = "Caja"
$ws2.Range("B7").Value = "Pagina web"
$ws2.Range("C7").Value = "url"

$ws2.Range("A8").Value = "Listado"
$ws2.Range("B8").Value = "Pais"
$ws2.Range("C8").Value = "numero"
$ws2.Range("D8").Value = "Selecciona un valor de la lista"
$ws2.Range("E8").Value = "si"
$ws2.Range("H8").Value = "y el pais?"
$ws2.Range("K8").Value = "Paises::nombre"

$ws2.Range("A9").Value = "Caja"
$ws2.Range("B9").Value = "Observaciones"
$ws2.Range("C9").Value = "mucho texto"
$ws2.Range("E9").Value = "no"
$ws2.Range("G9").Value = 100

$ws2.Range("A10").Value = "Caja"
$ws2.Range("B10").Value = "Fecha registro"
$ws2.Range("C10").Value = "fecha"
$ws2.Range("E10").Value = "si"

$ws2.Range("A11").Value = "Caja"
$ws2.Range("B11").Value = "Hora registro"
$ws2.Range("C11").Value = "hora"

$ws2.Range("A12").Value = "Caja"
$ws2.Range("B12").Value = "Ultima actualizacion"
$ws2.Range("C12").Value = "fecha hora"

$ws2.Range("A13").Value = "Cuadrado"
$ws2.Range("B13").Value = "Recibir publicidad"
$ws2.Range("C13").Value = "numero"
$ws2.Range("K13").Value = "1=Si"

# ======================================================================
# Sheet "ZeroCode": add the new "hora"/"fecha hora" data types and
# reorder the "Tipos de datos" list accordingly.
# ======================================================================
$ws3.Range("B2").Value = "texto"
$ws3.Range("B4").Value = "correo"
$ws3.Range("B5").Value = "clave"
$ws3.Range("B6").Value = "fecha"
$ws3.Range("B7").Value = "hora"
$ws3.Range("B8").Value = "fecha hora"
$ws3.Range("B9").Value = "mucho texto"
$ws3.Range("B10").Value = "url"

# ======================================================================
# Selections / active sheet & cell per sheet. "Clientes" must be
# selected/activated LAST so it ends up as the workbook's active tab.
# ======================================================================
$ws1.Range("E19").Select()
$ws3.Range("C15").Select()
$ws2.Range("C12").Select()
